# Auto-generated edit script: reorders/renames government-grants columns
# and fixes several jurisdiction/category labels across all 5 sheets,
# per the West Virginia overview workbook revision.

function Set-HeaderCell {
    param($ws, $r, $c, $val)
    $ws.Cells.Item($r, $c).Value = $val
}

function Set-TextCell {
    param($ws, $r, $c, $val)
    $ws.Cells.Item($r, $c).NumberFormat = "@"
    $ws.Cells.Item($r, $c).Value = $val
}

$wb = $excel.ActiveWorkbook

# ---- Sheet: Overall ----
$ws = $wb.Worksheets.Item('Overall')

$data = @(
    @('Share of 990 filers with government grants at risk', 'Number of 990 filers with government grants', 'Total government grants ($)', 'Size of operating surplus with government grants', 'Size of operating surplus without government grants'),
    @('75.81%', '769', '$1,067,242,685', '9.05%', '-27.48%')
)

for ($r = 1; $r -le $data.Length; $r++) {
    $row = $data[$r - 1]
    for ($c = 1; $c -le $row.Length; $c++) {
        if ($r -eq 1) {
            Set-HeaderCell $ws $r $c $row[$c - 1]
        } else {
            Set-TextCell $ws $r $c $row[$c - 1]
        }
    }
}

# ---- Sheet: County ----
$ws = $wb.Worksheets.Item('County')

$data = @(
    @('Geography', 'Share of 990 filers with government grants at risk', 'Number of 990 filers with government grants', 'Total government grants ($)', 'Size of operating surplus with government grants', 'Size of operating surplus without government grants'),
    @('United States', '67.35%', '103,475', '$267,700,640,005', '9.05%', '-12.83%'),
    @('West Virginia', '75.81%', '769', '$1,067,242,685', '9.05%', '-27.48%'),
    @('Barbour County', '75.00%', '8', '$9,024,036', '6.01%', '-11.96%'),
    @('Berkeley County', '75.00%', '28', '$19,567,722', '8.47%', '-11.77%'),
    @('Boone County', '100.00%', '6', '$5,253,491', '-1.85%', '-62.04%'),
    @('Braxton County', '50.00%', '4', '$2,612,747', '7.74%', '5.57%'),
    @('Brooke County', '83.33%', '6', '$9,836,958', '2.66%', '-24.27%'),
    @('Cabell County', '75.47%', '53', '$183,489,133', '3.68%', '-38.81%'),
    @('Calhoun County', '100.00%', '3', '$6,021,595', '15.16%', '-42.76%'),
    @('Doddridge County', '100.00%', '4', '$589,607', '-3.12%', '-85.30%'),
    @('Fayette County', '60.00%', '5', '$13,411,248', '32.58%', '-37.33%'),
    @('Gilmer County', '100.00%', '1', '$68,500', '-68.25%', '-94.96%'),
    @('Grant County', '100.00%', '5', '$3,221,163', '7.87%', '-14.97%'),
    @('Greenbrier County', '64.00%', '25', '$15,995,246', '19.24%', '-13.40%'),
    @('Hampshire County', '75.00%', '8', '$3,270,603', '9.15%', '-22.02%'),
    @('Hancock County', '76.92%', '13', '$13,711,243', '8.92%', '-31.85%'),
    @('Hardy County', '84.62%', '13', '$17,758,851', '4.58%', '-24.74%'),
    @('Harrison County', '87.50%', '24', '$24,506,375', '3.04%', '-31.20%'),
    @('Jackson County', '66.67%', '6', '$1,290,773', '6.61%', '-5.91%'),
    @('Jefferson County', '76.19%', '21', '$9,019,123', '19.50%', '-22.09%'),
    @('Kanawha County', '79.28%', '111', '$191,223,147', '6.21%', '-32.49%'),
    @('Lewis County', '90.00%', '10', '$3,836,222', '4.47%', '-31.16%'),
    @('Lincoln County', '83.33%', '6', '$8,473,846', '5.88%', '-70.12%'),
    @('Logan County', '84.62%', '13', '$10,105,029', '15.41%', '-38.87%'),
    @('Marion County', '50.00%', '22', '$25,630,232', '18.53%', '-2.99%'),
    @('Marshall County', '77.78%', '9', '$2,967,773', '9.47%', '-65.09%'),
    @('Mason County', '66.67%', '6', '$680,673', '5.17%', '-28.79%'),
    @('McDowell County', '100.00%', '7', '$7,019,129', '3.00%', '-85.58%'),
    @('Mercer County', '80.95%', '21', '$34,475,684', '8.14%', '-33.03%'),
    @('Mineral County', '78.57%', '14', '$9,860,443', '5.76%', '-30.83%'),
    @('Mingo County', '90.00%', '10', '$17,344,659', '13.54%', '-44.19%'),
    @('Monongalia County', '65.96%', '47', '$162,541,682', '15.88%', '-14.30%'),
    @('Monroe County', '75.00%', '4', '$719,401', '36.67%', '-36.46%'),
    @('Morgan County', '87.50%', '8', '$2,945,045', '17.69%', '-42.68%'),
    @('Nicholas County', '58.33%', '12', '$11,963,630', '13.26%', '-10.04%'),
    @('Ohio County', '68.18%', '44', '$47,587,191', '10.22%', '-11.97%'),
    @('Pendleton County', '75.00%', '8', '$4,709,394', '18.64%', '-21.02%'),
    @('Pocahontas County', '100.00%', '5', '$2,227,649', '7.19%', '-19.50%'),
    @('Preston County', '80.00%', '15', '$5,342,846', '11.99%', '-22.13%'),
    @('Putnam County', '55.56%', '9', '$8,740,115', '43.50%', '-13.32%'),
    @('Raleigh County', '74.07%', '27', '$42,781,625', '7.06%', '-36.57%'),
    @('Randolph County', '85.00%', '20', '$21,758,446', '6.12%', '-25.88%'),
    @('Ritchie County', '100.00%', '5', '$4,566,159', '17.17%', '-35.70%'),
    @('Roane County', '83.33%', '6', '$5,218,150', '-8.43%', '-22.50%'),
    @('Summers County', '75.00%', '4', '$2,184,887', '29.11%', '-60.63%'),
    @('Taylor County', '100.00%', '5', '$5,863,416', '6.98%', '-35.97%'),
    @('Tucker County', '85.71%', '7', '$7,766,231', '8.31%', '-25.75%'),
    @('Tyler County', '66.67%', '3', '$356,360', '29.44%', '-24.44%'),
    @('Upshur County', '58.33%', '12', '$18,081,957', '9.63%', '-5.97%'),
    @('Wayne County', '50.00%', '6', '$5,535,111', '20.36%', '-2.58%'),
    @('Webster County', '100.00%', '2', '$5,501,011', '-4.33%', '-25.82%'),
    @('Wetzel County', '55.56%', '9', '$3,848,787', '9.13%', '-5.83%'),
    @('Wirt County', '100.00%', '2', '$5,384,850', '1.31%', '-31.28%'),
    @('Wood County', '78.05%', '41', '$27,533,082', '11.19%', '-28.95%'),
    @('Wyoming County', '66.67%', '6', '$13,820,409', '7.01%', '-44.51%')
)

for ($r = 1; $r -le $data.Length; $r++) {
    $row = $data[$r - 1]
    for ($c = 1; $c -le $row.Length; $c++) {
        if ($r -eq 1) {
            Set-HeaderCell $ws $r $c $row[$c - 1]
        } else {
            Set-TextCell $ws $r $c $row[$c - 1]
        }
    }
}

# ---- Sheet: Congressional District ----
$ws = $wb.Worksheets.Item('Congressional District')

$data = @(
    @('Geography', 'Share of 990 filers with government grants at risk', 'Number of 990 filers with government grants', 'Total government grants ($)', 'Size of operating surplus with government grants', 'Size of operating surplus without government grants'),
    @('United States', '67.35%', '103,475', '$267,700,640,005', '9.05%', '-12.83%'),
    @('West Virginia', '75.81%', '769', '$1,067,242,685', '9.05%', '-27.48%'),
    @('Congressional District 1', '76.36%', '368', '$606,251,132', '8.07%', '-30.09%'),
    @('Congressional District 2', '75.31%', '401', '$460,991,553', '9.46%', '-24.08%')
)

for ($r = 1; $r -le $data.Length; $r++) {
    $row = $data[$r - 1]
    for ($c = 1; $c -le $row.Length; $c++) {
        if ($r -eq 1) {
            Set-HeaderCell $ws $r $c $row[$c - 1]
        } else {
            Set-TextCell $ws $r $c $row[$c - 1]
        }
    }
}

# ---- Sheet: Size ----
$ws = $wb.Worksheets.Item('Size')

$data = @(
    @('Size', 'Share of 990 filers with government grants at risk', 'Number of 990 filers with government grants', 'Total government grants ($)', 'Size of operating surplus with government grants', 'Size of operating surplus without government grants'),
    @('Between $100K and $499K', '75.90%', '278', '$39,348,992', '9.10%', '-35.34%'),
    @('Between $1M and $4.99M', '77.05%', '183', '$225,205,168', '9.29%', '-22.74%'),
    @('Between $500K and $999K', '75.68%', '111', '$42,122,126', '9.18%', '-32.13%'),
    @('Between $5M and $9.99M', '75.00%', '32', '$84,072,249', '11.32%', '-13.71%'),
    @('Greater than $10M', '73.40%', '94', '$671,068,805', '3.67%', '-8.91%'),
    @('Less than $100K', '76.06%', '71', '$5,425,345', '29.44%', '-32.25%'),
    @('Total', '75.81%', '769', '$1,067,242,685', '9.05%', '-27.48%')
)

for ($r = 1; $r -le $data.Length; $r++) {
    $row = $data[$r - 1]
    for ($c = 1; $c -le $row.Length; $c++) {
        if ($r -eq 1) {
            Set-HeaderCell $ws $r $c $row[$c - 1]
        } else {
            Set-TextCell $ws $r $c $row[$c - 1]
        }
    }
}

# ---- Sheet: Subsector ----
$ws = $wb.Worksheets.Item('Subsector')

$data = @(
    @('Subsector', 'Share of 990 filers with government grants at risk', 'Number of 990 filers with government grants', 'Total government grants ($)', 'Size of operating surplus with government grants', 'Size of operating surplus without government grants'),
    @('Arts, Culture, and Humanities', '71.05%', '38', '$16,481,528', '9.29%', '-20.83%'),
    @('Education (Excluding Universities)', '65.62%', '32', '$26,427,619', '15.46%', '-20.25%'),
    @('Environment and Animals', '51.61%', '31', '$7,559,689', '15.26%', '-0.53%'),
    @('Health (Excluding Hospitals)', '71.43%', '63', '$124,845,978', '8.92%', '-13.89%'),
    @('Hospitals', '61.54%', '13', '$36,156,605', '3.29%', '-2.04%'),
    @('Human Services', '79.75%', '316', '$187,409,162', '9.40%', '-37.74%'),
    @('Public, Societal Benefit', '81.82%', '55', '$136,546,982', '1.06%', '-42.53%'),
    @('Religion Related', '71.43%', '7', '$2,701,765', '2.96%', '-2.00%'),
    @('Unclassified', '76.70%', '206', '$510,666,753', '7.88%', '-19.49%'),
    @('Universities', '75.00%', '8', '$18,446,604', '2.76%', '-10.50%'),
    @('Total', '75.81%', '769', '$1,067,242,685', '9.05%', '-27.48%')
)

for ($r = 1; $r -le $data.Length; $r++) {
    $row = $data[$r - 1]
    for ($c = 1; $c -le $row.Length; $c++) {
        if ($r -eq 1) {
            Set-HeaderCell $ws $r $c $row[$c - 1]
        } else {
            Set-TextCell $ws $r $c $row[$c - 1]
        }
    }
}

